$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Clojure" term: drop the spellStart/spellEnd proofErr markers that wrap
#    the run, while keeping the bold "Clojure" run + the plain description
#    run (and the list-paragraph / numbering properties) exactly as-is.
#    proofErr markers have no text width, so Find/Replace can't touch them -
#    we replace the whole paragraph's XML instead.
# ---------------------------------------------------------------------------
$clojureXml = @'
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:b/></w:rPr><w:t>Clojure</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/></w:rPr><w:t>- Dialect of LISP (second oldest programming language, highly complex), general-purpose program with an emphasis on functional programing (style of building the structure and elements of a computer program).</w:t></w:r></w:p>
'@

$rClojure = $d.Content
$rClojure.Find.Execute("Clojure- Dialect of LISP") | Out-Null
$pClojure = $rClojure.Paragraphs(1)
$pClojure.Range.InsertXML($clojureXml)

# ---------------------------------------------------------------------------
# 2) Remove the heading block that currently sits right after the
#    "User Story" paragraph (it gets re-created, with new content mixed in,
#    earlier in the flow - see step 3). Locate it by its start/end text so
#    we don't depend on fragile paragraph indices.
# ---------------------------------------------------------------------------
$rOldStart = $d.Content
$rOldStart.Find.Execute("What's With All These Conferences, Anyway?") | Out-Null
$oldBlockStart = $rOldStart.Paragraphs(1).Range.Start

$rOldEnd = $d.Content
$rOldEnd.Find.Execute("Look How Big and Weird Things Get") | Out-Null
$oldBlockEnd = $rOldEnd.Paragraphs(1).Range.End

$oldBlockRange = $d.Range($oldBlockStart, $oldBlockEnd)
$oldBlockRange.Delete()

# ---------------------------------------------------------------------------
# 3) Right after "User Story: (use case) ... users." :
#      - drop the bookmarkStart/bookmarkEnd that currently lives at the end
#        of that paragraph (it moves into its own new empty paragraph below)
#      - insert the (moved) heading block
#      - insert four new term definitions: Python, Abstraction, Jython,
#        IronPython
#      - insert a new, otherwise-empty paragraph that now owns the
#        _GoBack bookmark
# ---------------------------------------------------------------------------
$userStoryXml = @'
<w:p><w:pPr><w:rPr><w:color w:val="222222"/><w:highlight w:val="white"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:color w:val="222222"/><w:highlight w:val="white"/></w:rPr><w:t xml:space="preserve">User Story: </w:t></w:r><w:r><w:rPr><w:color w:val="222222"/><w:highlight w:val="white"/></w:rPr><w:t xml:space="preserve">(use case) what the program will look like to the users. </w:t></w:r></w:p>
'@

$rUserStory = $d.Content
$rUserStory.Find.Execute("User Story:") | Out-Null
$pUserStory = $rUserStory.Paragraphs(1)
$pUserStory.Range.InsertXML($userStoryXml)

$newBlockXml = @'
<w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">What's </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>With</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> All These Conferences, Anyway?</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Why Are Programmers So Intense About Languages?</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>The Beauty of the Standard Library</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>What Do Different Languages Do?</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>The Importance of C</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>The Corporate Object Revolution</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Look How Big and Weird Things Get </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>With</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> Just Python</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Python</w:t></w:r><w:r><w:t>: Object oriented language that is built on many abstractions to work with different languages, such as C, C++ and Fortran77</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Abstraction</w:t></w:r><w:r><w:t xml:space="preserve">: Hiding all but the relevant data about an object </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>in order to</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> reduce complexity and increase efficiency.</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>Jython</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: Version of Python designed to run inside of Java</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>IronPython</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: Version of Python that works with Microsoft\u2019s .NET</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$newBlockXml = $newBlockXml.Replace("\u2019", [char]0x2019)

$rUserStory2 = $d.Content
$rUserStory2.Find.Execute("User Story:") | Out-Null
$pUserStory2 = $rUserStory2.Paragraphs(1)
$insertPoint = $d.Range($pUserStory2.Range.End, $pUserStory2.Range.End)
$insertPoint.InsertXML($newBlockXml)

# ---------------------------------------------------------------------------
# 4) Drop the stray <w:lastRenderedPageBreak/> that used to sit in front of
#    "Why Are There So Many Languages?" (the page-break hint now effectively
#    belongs to the new IronPython paragraph instead).
# ---------------------------------------------------------------------------
$whyThereXml = @'
<w:p><w:pPr><w:spacing w:before="160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Why Are There So Many Languages?</w:t></w:r></w:p>
'@

$rWhyThere = $d.Content
$rWhyThere.Find.Execute("Why Are There So Many Languages?") | Out-Null
$pWhyThere = $rWhyThere.Paragraphs(1)
$pWhyThere.Range.InsertXML($whyThereXml)

Write-Output "done"
